# Swap the data (columns B:AC) between row pairs (11,12) and (143,144).
# Column A (the rank/id number) stays with its original row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowPairs = @(
    @(11, 12),
    @(143, 144)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $rng1 = $ws.Range("B" + $r1 + ":AC" + $r1)
    $rng2 = $ws.Range("B" + $r2 + ":AC" + $r2)

    $vals1 = $rng1.Value2
    $vals2 = $rng2.Value2

    $rng1.Value2 = $vals2
    $rng2.Value2 = $vals1
}
